# Update results values on each year sheet (row 2) with the latest
# server-computed results, per the "ADD results from server" commit.

$wb = $excel.ActiveWorkbook

# Sheet "2025" -> its own unique set of updated values
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 545.7008988199987
$ws2025.Range("E2").Value = 21271.60964344695
$ws2025.Range("I2").Value = 10981.098647904
$ws2025.Range("L2").Value = 42839.69549276341
$ws2025.Range("M2").Value = 8019.132604175002
$ws2025.Range("N2").Value = 4536.424898579853
$ws2025.Range("O2").Value = 5082.110035902254

# Sheets "2030", "2035", "2040", "2045", "2050" -> all share the same
# updated values as one another
$years = @("2030", "2035", "2040", "2045", "2050")
foreach ($year in $years) {
    $ws = $wb.Worksheets.Item($year)
    $ws.Range("A2").Value = 883.0954041229934
    $ws.Range("B2").Value = 4582.57765438246
    $ws.Range("E2").Value = 50104.2345376395
    $ws.Range("I2").Value = 37467.9321740052
    $ws.Range("L2").Value = 64560.9394462146
    $ws.Range("M2").Value = 20366.03949924198
    $ws.Range("N2").Value = 11003.6745215972
    $ws.Range("O2").Value = 10532.83682572548
}
